$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1831.591157525936
$ws.Range("B3").Value = 2900.686774403619
$ws.Range("B4").Value = 3895.74006538886
$ws.Range("B5").Value = 4359.462888536614
$ws.Range("B6").Value = 4908.677389009084
$ws.Range("B7").Value = 5486.642319586107
$ws.Range("B8").Value = 5700.138037599012
$ws.Range("B9").Value = 6001.457960867561
$ws.Range("B10").Value = 6347.371853381838
$ws.Range("B11").Value = 6511.84542830972
$ws.Range("B12").Value = 6794.667971583616
$ws.Range("B13").Value = 6987.925687287144
$ws.Range("B14").Value = 7074.915236817072
$ws.Range("B15").Value = 7302.369183577823
$ws.Range("B16").Value = 7425.296298002776
$ws.Range("B17").Value = 7639.665539629861
$ws.Range("B18").Value = 7705.27926169254
$ws.Range("B19").Value = 7895.640957707441
$ws.Range("B20").Value = 7998.670052572564
$ws.Range("B21").Value = 8037.150975487283
$ws.Range("B22").Value = 8150.445646074456
$ws.Range("B23").Value = 8241.785704796668
$ws.Range("B24").Value = 8273.790333035993
$ws.Range("B25").Value = 8280.59494368448
$ws.Range("B26").Value = 8292.317813205944
$ws.Range("B27").Value = 8293.00593833287
$ws.Range("B28").Value = 8293.00593833287
$ws.Range("B29").Value = 8293.00593833287
$ws.Range("B30").Value = 8293.00593833287
$ws.Range("B31").Value = 8293.00593833287
$ws.Range("B32").Value = 8293.00593833287
$ws.Range("B33").Value = 8293.00593833287
$ws.Range("B34").Value = 8293.00593833287
$ws.Range("B35").Value = 8293.00593833287
$ws.Range("B36").Value = 8293.00593833287
$ws.Range("B37").Value = 8293.00593833287
$ws.Range("B38").Value = 8293.00593833287
$ws.Range("B39").Value = 8293.00593833287
$ws.Range("B40").Value = 8293.00593833287
$ws.Range("B41").Value = 8293.00593833287
$ws.Range("B42").Value = 8293.00593833287
$ws.Range("B43").Value = 8293.00593833287
$ws.Range("B44").Value = 8293.00593833287
$ws.Range("B45").Value = 8293.00593833287
$ws.Range("B46").Value = 8293.00593833287
$ws.Range("B47").Value = 8293.00593833287
$ws.Range("B48").Value = 8293.00593833287
$ws.Range("B49").Value = 8293.00593833287
$ws.Range("B50").Value = 8293.00593833287
$ws.Range("B51").Value = 8293.00593833287
$ws.Range("B52").Value = 8293.00593833287
$ws.Range("B53").Value = 8293.00593833287
$ws.Range("B54").Value = 8293.00593833287
$ws.Range("B55").Value = 8293.00593833287
$ws.Range("B56").Value = 8293.00593833287
$ws.Range("B57").Value = 8293.00593833287
$ws.Range("B58").Value = 8293.00593833287
$ws.Range("B59").Value = 8293.00593833287
$ws.Range("B60").Value = 8293.00593833287
$ws.Range("B61").Value = 8293.00593833287
$ws.Range("B62").Value = 8293.00593833287
